$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 LogisticRegression(class_weight=''balanced'', l1_ratio=0.5,
                                    max_iter=1000, penalty=''elasticnet'',
                                    random_state=42, solver=''saga''))])'
$ws.Range("B2").Value = 0.6445565545565545
$ws.Range("C2").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.5, ''model__class_weight'': ''balanced''}'
$ws.Range("D2").Value = 0.6731902416603736
$ws.Range("E2").Value = 0.6055782412032412
$ws.Range("F2").Value = 0.75
$ws.Range("G2").Value = 0.7353542007662698
$ws.Range("H2").Value = 0.6564484126984128
$ws.Range("I2").Value = 0.75
$ws.Range("J2").Value = 0.6220744680851064
$ws.Range("K2").Value = 0.5888888888888889
$ws.Range("L2").Value = 0.75
$ws.Range("N2").Value = '[1 0 1 1 1 1 0 0 1 0 0 1 1 1 0 0 1 1 1 1 1 1 1 0]'

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f17783deee0>),
                (''model'',
                 LogisticRegression(class_weight=''balanced'', l1_ratio=0.1,
                                    max_iter=1000, penalty=''elasticnet'',
                                    random_state=42, solver=''saga''))])'
$ws.Range("B3").Value = 0.5828282828282828
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f17783dc2e0>, ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.1, ''model__class_weight'': ''balanced''}'
$ws.Range("D3").Value = 0.6682811777152707
$ws.Range("E3").Value = 0.4878713416213416
$ws.Range("F3").Value = 0.7999999999999999
$ws.Range("G3").Value = 0.7280682840970768
$ws.Range("H3").Value = 0.6339583333333333
$ws.Range("I3").Value = 0.7368421052631579
$ws.Range("J3").Value = 0.6197695035460993
$ws.Range("K3").Value = 0.4159722222222222
$ws.Range("L3").Value = 0.875
$ws.Range("N3").Value = '[1 1 1 1 1 0 1 1 1 1 1 0 1 0 1 1 1 1 1 0 1 0 1 1]'

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f17783dcfd0>),
                (''model'',
                 LogisticRegression(class_weight=''balanced'', l1_ratio=0.01,
                                    max_iter=1000, penalty=''elasticnet'',
                                    random_state=42, solver=''saga''))])'
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f177841bac0>, ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.01, ''model__class_weight'': ''balanced''}'
$ws.Range("D4").Value = 0.658301352095244
$ws.Range("E4").Value = 0.561356721981722
$ws.Range("G4").Value = 0.7221039912805383
$ws.Range("H4").Value = 0.6739484126984125
$ws.Range("J4").Value = 0.6057407407407407
$ws.Range("K4").Value = 0.5091666666666667

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 LogisticRegression(class_weight=''balanced'', l1_ratio=0.5,
                                    max_iter=1000, penalty=''elasticnet'',
                                    random_state=42, solver=''saga''))])'
$ws.Range("B5").Value = 0.6424741924741924
$ws.Range("C5").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.5, ''model__class_weight'': ''balanced''}'
$ws.Range("D5").Value = 0.6769280689154774
$ws.Range("E5").Value = 0.5845659664409664
$ws.Range("F5").Value = 0.6153846153846153
$ws.Range("G5").Value = 0.735641343965263
$ws.Range("H5").Value = 0.6379828042328042
$ws.Range("I5").Value = 0.6666666666666666
$ws.Range("J5").Value = 0.6292517006802721
$ws.Range("K5").Value = 0.5743055555555555
$ws.Range("L5").Value = 0.5714285714285714
$ws.Range("N5").Value = '[0 0 1 1 1 0 0 0 0 0 1 1 1 1 1 0 0 0 1 0 0 1 1 1]'

# Row 6
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f177841b460>),
                (''model'',
                 LogisticRegression(class_weight=''balanced'', l1_ratio=0.5,
                                    max_iter=1000, penalty=''elasticnet'',
                                    random_state=42, solver=''saga''))])'
$ws.Range("B6").Value = 0.7176845376845377
$ws.Range("C6").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f17782eee80>, ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.5, ''model__class_weight'': ''balanced''}'
$ws.Range("D6").Value = 0.7553392646378311
$ws.Range("E6").Value = 0.6077489639989641
$ws.Range("F6").Value = 0.5925925925925926
$ws.Range("G6").Value = 0.7992066077896293
$ws.Range("H6").Value = 0.6583779761904762
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0.7174679487179487
$ws.Range("K6").Value = 0.5861111111111111
$ws.Range("L6").Value = 0.7272727272727273
$ws.Range("N6").Value = '[1 1 1 0 1 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 1 1 1 1]'
